$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.233.16'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.375.70'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.695'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.26'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +5.37%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +24.52%  '
$ws.Range('E10').Value = '  +5.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.00'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.29'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +19.28%  '
$ws.Range('E13').Value = '  +19.05%  '
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').Value = '2.729.41'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('E16').Value = '  +6.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.923'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +7.66%  '
$ws.Range('D18').Value = '2.370.73'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').Value = '44.342.45'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('E20').Value = '  +4.76%  '
$ws.Range('E21').Value = '  +5.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.60'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '258.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.70'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +8.14%  '
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('E30').Value = '  +2.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.37'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('E33').Value = '  +6.82%  '
$ws.Range('E34').Value = '  +8.04%  '
$ws.Range('E35').Value = '  +10.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.35'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.50%  '
$ws.Range('E37').Value = '  +6.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.49'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('E40').Value = '  +7.83%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.15'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.19'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.22'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.194'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +16.05%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.101'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +5.75%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.27'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.52'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +12.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.40'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '1.477.75'
$ws.Range('E51').Value = '  +1.82%  '
